$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7769267559051514
$ws.Range("B1").Value = 3.181809186935425
$ws.Range("C1").Value = 2.912545204162598
$ws.Range("D1").Value = 2.468739032745361
$ws.Range("E1").Value = 2.135701417922974
